# [PV-350][WIP] Replace hard coding of visual height with calculated value
#
# Update the header row labels on the main plan sheet and move the
# active cell selection from F1 to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"

$ws.Range("F2").Select()
